$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "H 72" record (row 2). All rows below shift up by one,
# shrinking the used range from A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()
